$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new data rows (A6:B6 and A7:B7), preserving the existing
# pattern of text-typed "ID" values (e.g. "0005", "0006") and the
# repeated "CO" type value, without introducing any new cell styles. ---

# Build the text values via a formula first (so the strings keep their
# leading zeros and are not auto-coerced to numbers), then convert the
# formulas to static values in-place using Copy / PasteSpecial(values).
$ws.Range("A6").Formula = '="0005"'
$ws.Range("A7").Formula = '="0006"'
$ws.Range("A6:A7").Copy()
$ws.Range("A6:A7").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B6").Value = "CO"
$ws.Range("B7").Value = "CO"

# Touch D6/D7 so an (empty) cell entry is materialized for those rows,
# matching the existing D2:D5 pattern.
$ws.Range("D6").Borders.LineStyle = 0
$ws.Range("D7").Borders.LineStyle = 0

# --- Add cell comments (legacy/"Leyenda" notes) on C1 and E1 ---
$c1 = $ws.Range("C1").AddComment("test 1")
$c2 = $ws.Range("E1").AddComment("test 3")
